$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.922.07"
$ws.Range("E2").Value = "  -1.43%  "

$ws.Range("D3").Value = "2.225.04"
$ws.Range("E3").Value = "  -0.83%  "

$ws.Range("E4").Value = "  -1.85%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "90.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.553"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.66%  "

$ws.Range("E8").Value = "  -0.64%  "

$ws.Range("E9").Value = "  -6.39%  "

$ws.Range("E10").Value = "  -5.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0777"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.07%  "

$ws.Range("D14").Value = "2.564.46"
$ws.Range("E14").Value = "  -0.85%  "

$ws.Range("D15").Value = "2.228.89"
$ws.Range("E15").Value = "  -2.60%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.775"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.15%  "

$ws.Range("D18").Value = "43.761.78"
$ws.Range("E18").Value = "  -0.96%  "

$ws.Range("D19").Value = "0.0₃0901"
$ws.Range("E19").Value = "  -5.89%  "

$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.82%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.84%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.43%  "

$ws.Range("E25").Value = "  -0.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.80%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "38.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.96%  "

$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "152.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.97%  "

$ws.Range("E32").Value = "  -9.38%  "

$ws.Range("E33").Value = "  -5.77%  "

$ws.Range("E34").Value = "  -5.26%  "

$ws.Range("E35").Value = "  -3.73%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.102"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.26%  "

$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.70%  "

$ws.Range("E38").Value = "  -5.09%  "

$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.29%  "

$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").Value = "1.838.36"
$ws.Range("E44").Value = "  +4.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.182"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "67.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.49%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "94.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.11%  "

$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "73.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "13.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.09%  "

$ws.Range("D51").Value = "2.446.01"
$ws.Range("E51").Value = "  -0.78%  "

Write-Host "cryptos list updated"
